$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44875
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = 1550
$ws.Range("P2").Value = 1550

# Row 3
$ws.Range("D3").Value = 45245
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1900
$ws.Range("P3").Value = 1900

# Row 4
$ws.Range("D4").Value = 44510
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 1300
$ws.Range("L4").Value = 1400
$ws.Range("M4").Value = 1350
$ws.Range("N4").Value = '$/kilo'
$ws.Range("P4").Value = 1350

# Row 5
$ws.Range("D5").Value = 45203
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 1400
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1450
$ws.Range("O5").Value = 'Provincia de Linares'
$ws.Range("P5").Value = 1450

# Row 6
$ws.Range("D6").Value = 44468
$ws.Range("H6").Value = 'Verde'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 1800
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 1920
$ws.Range("O6").Value = 'Provincia de Linares'
$ws.Range("P6").Value = 1920

# Row 7
$ws.Range("D7").Value = 45202
$ws.Range("H7").Value = 'Verde'
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 1600
$ws.Range("L7").Value = 1600
$ws.Range("M7").Value = 1600
$ws.Range("P7").Value = 1600

# Row 8
$ws.Range("D8").Value = 44496
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 550
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 1773
$ws.Range("N8").Value = '$/paquete'
$ws.Range("P8").Value = 1773

# Row 9
$ws.Range("D9").Value = 45246
$ws.Range("K9").Value = 2000
$ws.Range("M9").Value = 2000
$ws.Range("P9").Value = 2000

# Row 10
$ws.Range("D10").Value = 45246
$ws.Range("I10").Value = 'Segunda'
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = 1500
$ws.Range("P10").Value = 1500

# Row 11
$ws.Range("D11").Value = 44489
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 1400
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = 1450
$ws.Range("P11").Value = 1450

# Row 12
$ws.Range("D12").Value = 44524
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 1500
$ws.Range("L12").Value = 1600
$ws.Range("M12").Value = 1550
$ws.Range("O12").Value = 'Provincia de Talca'
$ws.Range("P12").Value = 1550

# Row 13
$ws.Range("D13").Value = 44839
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 1700
$ws.Range("L13").Value = 1800
$ws.Range("M13").Value = 1760
$ws.Range("P13").Value = 1760

# Row 14
$ws.Range("D14").Value = 44545
$ws.Range("H14").Value = 'Sin especificar'
$ws.Range("J14").Value = 550
$ws.Range("K14").Value = 1700
$ws.Range("L14").Value = 1800
$ws.Range("M14").Value = 1755
$ws.Range("P14").Value = 1755

# Row 15
$ws.Range("D15").Value = 44881
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 2600
$ws.Range("L15").Value = 2700
$ws.Range("M15").Value = 2650
$ws.Range("P15").Value = 2650

# Row 16
$ws.Range("D16").Value = 44881
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 2400
$ws.Range("L16").Value = 2400
$ws.Range("M16").Value = 2400
$ws.Range("P16").Value = 2400

# Row 17
$ws.Range("D17").Value = 44868
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 1300
$ws.Range("M17").Value = 1250
$ws.Range("O17").Value = 'Región del Maule'
$ws.Range("P17").Value = 1250

# Row 18
$ws.Range("D18").Value = 44868
$ws.Range("I18").Value = 'Segunda'
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 1000
$ws.Range("O18").Value = 'Región del Maule'
$ws.Range("P18").Value = 1000

# Row 19
$ws.Range("D19").Value = 44519
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 1200
$ws.Range("L19").Value = 1300
$ws.Range("M19").Value = 1240
$ws.Range("P19").Value = 1240

# Row 20
$ws.Range("D20").Value = 44876
$ws.Range("J20").Value = 350
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 1600
$ws.Range("M20").Value = 1557
$ws.Range("P20").Value = 1557

# Row 21
$ws.Range("D21").Value = 45230
$ws.Range("J21").Value = 500
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 1500
$ws.Range("P21").Value = 1500

# Row 22
$ws.Range("D22").Value = 44511
$ws.Range("H22").Value = 'Sin especificar'
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 1300
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = 1350
$ws.Range("P22").Value = 1350

# Row 23
$ws.Range("D23").Value = 44526
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 1600
$ws.Range("M23").Value = 1550
$ws.Range("P23").Value = 1550

# Row 24
$ws.Range("D24").Value = 44860
$ws.Range("J24").Value = 1100
$ws.Range("L24").Value = 1700
$ws.Range("M24").Value = 1609
$ws.Range("P24").Value = 1609

# Row 25
$ws.Range("D25").Value = 44477
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 1400
$ws.Range("L25").Value = 1500
$ws.Range("M25").Value = 1460
$ws.Range("O25").Value = 'Provincia de Linares'
$ws.Range("P25").Value = 1460
